$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the newly-logged meeting (row 23 / entry #21)
# Note: shared-string table order follows first-use order, so register
# "regression test" (F23) before "2/12" (B23) to match the original authoring order.
$ws.Range("F23").Value = "regression test"
$ws.Range("C23").Value = "7-11pm"
$ws.Range("D23").Value = "library317"
$ws.Range("E23").Value = "everyone"
$ws.Range("B23").Value = "2/12"

# Update the active selection to reflect where the user left off editing
$ws.Activate()
$ws.Range("B23").Select()
